$d = $word.ActiveDocument

# --- Sentencing-summary table -------------------------------------------
# The "Jail Days" row and the "Jail Days Suspended" row each report a
# numeric day count for both charges. The defendant only ever received
# jail-time *credit* (see below) and no jail time was actually imposed or
# suspended, so every numeric entry in those two rows becomes "None".
$t = $d.Tables.Item(1)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $label = $t.Cell($r, 1).Range.Text
    $label = $label.Replace([char]13, "").Replace([char]7, "")
    if ($label -eq "Jail Days" -or $label -eq "Jail Days Suspended") {
        for ($c = 2; $c -le $t.Columns.Count; $c++) {
            $t.Cell($r, $c).Range.Text = "None"
        }
    }
}

# --- Jail-time-credit sentence ------------------------------------------
# "...Defendant shall receive credit for 1 days already served in jail."
# becomes "...shall receive credit for 5 days already served in jail."
#
# The digit and the trailing space live in their own run, immediately
# followed by a separate (identically-formatted) run containing "days".
# A plain Range.Text assignment would cause the two adjacent same-format
# runs to coalesce into one, so the original run's FormattedText is
# captured first and re-applied after the text swap to keep the run
# boundary (and its own w:r/w:rPr) intact, matching the original markup
# shape as closely as possible.
$rng = $d.Content
$found = $rng.Find.Execute("credit for 1 days", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)
if ($found) {
    $prefixLen = "credit for ".Length
    $digitStart = $rng.Start + $prefixLen

    $runRange = $d.Range($digitStart, $digitStart + 2)   # "1 " (digit + space)
    $template = $runRange.FormattedText

    $digitOnly = $d.Range($digitStart, $digitStart + 1)
    $digitOnly.Text = "5"

    $runRange2 = $d.Range($digitStart, $digitStart + 2)  # "5 " after the edit
    $runRange2.FormattedText = $template
}

Write-Output ("Jail days credit sentence updated: " + $found)
